$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "58.256.45"
$ws.Range("E2").Value = "  -0.69%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.595.61"
$ws.Range("E3").Value = "  -1.20%  "
$ws.Range("E4").Value = "  +0.17%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "521.56"
$ws.Range("E5").Value = "  -0.56%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "143.50"
$ws.Range("E6").Value = "  +0.10%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.998"
$ws.Range("E7").Value = "  -0.05%  "
$ws.Range("E8").Value = "  -0.32%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.621.10"
$ws.Range("E9").Value = "  -0.61%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.48"
$ws.Range("E10").Value = "  -1.56%  "
$ws.Range("E11").Value = "  -1.78%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.340"
$ws.Range("E12").Value = "  +1.20%  "
$ws.Range("E13").Value = "  -0.09%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.056.62"
$ws.Range("E14").Value = "  -1.07%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "58.226.67"
$ws.Range("E15").Value = "  -0.61%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "20.37"
$ws.Range("E16").Value = "  -3.34%  "
$ws.Range("B17").Value = "ShibaInu"
$ws.Range("C17").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000134"
$ws.Range("E17").Value = "  -1.85%  "
$ws.Range("B18").Value = "WrappedEther"
$ws.Range("C18").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.606.34"
$ws.Range("E18").Value = "  +0.49%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "339.22"
$ws.Range("E19").Value = "  +0.21%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.32"
$ws.Range("E20").Value = "  -2.21%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.23"
$ws.Range("E21").Value = "  -2.08%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.44"
$ws.Range("E22").Value = "  +2.35%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.999"
$ws.Range("E23").Value = "  -0.09%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "65.44"
$ws.Range("E24").Value = "  +0.40%  "
$ws.Range("E25").Value = "  +0.71%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.404"
$ws.Range("E26").Value = "  -2.92%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.722.49"
$ws.Range("E27").Value = "  -0.41%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.996"
$ws.Range("E28").Value = "  -0.08%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.05"
$ws.Range("E29").Value = "  -1.51%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0₃0743"
$ws.Range("E30").Value = "  -6.59%  "
$ws.Range("E31").Value = "  -0.07%  "
$ws.Range("E32").Value = "  -6.11%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.58"
$ws.Range("E33").Value = "  -1.19%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "18.77"
$ws.Range("E34").Value = "  -0.34%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "149.17"
$ws.Range("E35").Value = "  -0.78%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.00"
$ws.Range("E36").Value = "  -3.09%  "
$ws.Range("E37").Value = "  -4.72%  "
$ws.Range("E38").Value = "  -5.18%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.857"
$ws.Range("E39").Value = "  +0.11%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.46"
$ws.Range("E40").Value = "  +1.67%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "36.09"
$ws.Range("E41").Value = "  -1.07%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.53"
$ws.Range("E42").Value = "  -2.86%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.997"
$ws.Range("E43").Value = "  -0.14%  "
$ws.Range("E44").Value = "  +0.68%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "271.02"
$ws.Range("E45").Value = "  +0.28%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "10.68"
$ws.Range("E46").Value = "  +0.24%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0955"
$ws.Range("E47").Value = "  -1.69%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "18.76"
$ws.Range("E48").Value = "  -1.99%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0521"
$ws.Range("E49").Value = "  -2.34%  "
$ws.Range("B50").Value = "Maker"
$ws.Range("C50").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.967.72"
$ws.Range("E50").Value = "  -3.30%  "
$ws.Range("B51").Value = "RenderToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "4.68"
$ws.Range("E51").Value = "  +1.90%  "
